$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("A9").Value = 112171798
$ws.Range("B9").Value = 78713
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 6458
$ws.Range("F9").Value = "Lunglav"
$ws.Range("G9").Value = "Lobaria pulmonaria"
$ws.Range("H9").Value = "(L.) Hoffm."
$ws.Range("Q9").Value = 756371
$ws.Range("R9").Value = 7212116
$ws.Range("AJ9").ClearContents()
$ws.Range("AK9").ClearContents()
$ws.Range("AO9").ClearContents()

# Row 10
$ws.Range("A10").Value = 112171776
$ws.Range("B10").Value = 85448
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 3739
$ws.Range("F10").Value = "Persiljespindling"
$ws.Range("G10").Value = "Cortinarius sulfurinus"
$ws.Range("H10").Value = "Quél."
$ws.Range("Q10").Value = 756261
$ws.Range("R10").Value = 7211953
$ws.Range("AJ10").ClearContents()
$ws.Range("AK10").ClearContents()
$ws.Range("AO10").ClearContents()

# Row 11
$ws.Range("A11").Value = 112171779
$ws.Range("B11").Value = 102192
$ws.Range("D11").Value = "LC"
$ws.Range("E11").Value = 222412
$ws.Range("F11").Value = "Tibast"
$ws.Range("G11").Value = "Daphne mezereum"
$ws.Range("H11").Value = "L."
$ws.Range("Q11").Value = 756291
$ws.Range("R11").Value = 7211892
$ws.Range("AJ11").ClearContents()
$ws.Range("AK11").ClearContents()
$ws.Range("AO11").ClearContents()

# Row 12
$ws.Range("A12").Value = 112171785
$ws.Range("B12").Value = 78677
$ws.Range("D12").Value = "LC"
$ws.Range("E12").Value = 229748
$ws.Range("F12").Value = "Gytterlav"
$ws.Range("G12").Value = "Protopannaria pezizoides"
$ws.Range("H12").Value = "(Weber) P.M.Jørg. & S.Ekman"
$ws.Range("Q12").Value = 756412
$ws.Range("R12").Value = 7211954
$ws.Range("AJ12").ClearContents()
$ws.Range("AK12").ClearContents()
$ws.Range("AO12").ClearContents()

# Row 13
$ws.Range("A13").Value = 112171792
$ws.Range("B13").Value = 85401
$ws.Range("D13").Value = "LC"
$ws.Range("E13").Value = 249228
$ws.Range("F13").Value = "Barrfagerspindling"
$ws.Range("G13").Value = "Cortinarius piceae"
$ws.Range("H13").Value = "Frøslev, T.S.Jeppesen & Brandrud"
$ws.Range("Q13").Value = 756395
$ws.Range("R13").Value = 7211974
$ws.Range("AJ13").ClearContents()
$ws.Range("AK13").ClearContents()
$ws.Range("AO13").ClearContents()

# Row 14
$ws.Range("A14").Value = 112171814
$ws.Range("B14").Value = 89499
$ws.Range("D14").Value = "NT"
$ws.Range("E14").Value = 112
$ws.Range("F14").Value = "Stjärntagging"
$ws.Range("G14").Value = "Asterodon ferruginosus"
$ws.Range("H14").Value = "Pat."
$ws.Range("Q14").Value = 756486
$ws.Range("R14").Value = 7212041
$ws.Range("AJ14").Value = "gran"
$ws.Range("AK14").Value = "Picea abies"
$ws.Range("AO14").Value = "Picea abies"

# Row 15
$ws.Range("A15").Value = 112171812
$ws.Range("B15").Value = 78746
$ws.Range("D15").Value = "LC"
$ws.Range("E15").Value = 6463
$ws.Range("F15").Value = "Bårdlav"
$ws.Range("G15").Value = "Nephroma parile"
$ws.Range("H15").Value = "(Ach.) Ach."
$ws.Range("Q15").Value = 756485
$ws.Range("R15").Value = 7212023
$ws.Range("AJ15").Value = "sälg"
$ws.Range("AK15").Value = "Salix caprea"
$ws.Range("AO15").Value = "Salix caprea"

# Row 16
$ws.Range("A16").Value = 112171787
$ws.Range("B16").Value = 86371
$ws.Range("D16").Value = "NT"
$ws.Range("E16").Value = 4412
$ws.Range("F16").Value = "Äggvaxskivling"
$ws.Range("G16").Value = "Hygrophorus karstenii"
$ws.Range("H16").Value = "Sacc. & Cub."
$ws.Range("Q16").Value = 756408
$ws.Range("R16").Value = 7211956
$ws.Range("AJ16").ClearContents()
$ws.Range("AK16").ClearContents()
$ws.Range("AO16").ClearContents()

# Row 17
$ws.Range("A17").Value = 112171806
$ws.Range("B17").Value = 86371
$ws.Range("D17").Value = "NT"
$ws.Range("E17").Value = 4412
$ws.Range("F17").Value = "Äggvaxskivling"
$ws.Range("G17").Value = "Hygrophorus karstenii"
$ws.Range("H17").Value = "Sacc. & Cub."
$ws.Range("Q17").Value = 756477
$ws.Range("R17").Value = 7212031
$ws.Range("AJ17").ClearContents()
$ws.Range("AK17").ClearContents()
$ws.Range("AO17").ClearContents()

# Row 18
$ws.Range("A18").Value = 112171810
$ws.Range("B18").Value = 90480
$ws.Range("D18").Value = "LC"
$ws.Range("E18").Value = 4769
$ws.Range("F18").Value = "Svavelriska"
$ws.Range("G18").Value = "Lactarius scrobiculatus"
$ws.Range("H18").Value = "(Scop.:Fr.) Fr."
$ws.Range("Q18").Value = 756486
$ws.Range("R18").Value = 7212020
$ws.Range("AJ18").ClearContents()
$ws.Range("AK18").ClearContents()
$ws.Range("AO18").ClearContents()

# Row 19
$ws.Range("A19").Value = 112171801
$ws.Range("B19").Value = 78713
$ws.Range("D19").Value = "NT"
$ws.Range("E19").Value = 6458
$ws.Range("F19").Value = "Lunglav"
$ws.Range("G19").Value = "Lobaria pulmonaria"
$ws.Range("H19").Value = "(L.) Hoffm."
$ws.Range("Q19").Value = 756448
$ws.Range("R19").Value = 7212052
$ws.Range("AJ19").ClearContents()
$ws.Range("AK19").ClearContents()
$ws.Range("AO19").ClearContents()

# Row 20
$ws.Range("A20").Value = 112171788
$ws.Range("B20").Value = 78713
$ws.Range("D20").Value = "NT"
$ws.Range("E20").Value = 6458
$ws.Range("F20").Value = "Lunglav"
$ws.Range("G20").Value = "Lobaria pulmonaria"
$ws.Range("H20").Value = "(L.) Hoffm."
$ws.Range("Q20").Value = 756401
$ws.Range("R20").Value = 7211954
$ws.Range("AJ20").Value = "sälg"
$ws.Range("AK20").Value = "Salix caprea"
$ws.Range("AO20").Value = "Salix caprea"

# Row 21
$ws.Range("A21").Value = 112171813
$ws.Range("B21").Value = 78713
$ws.Range("D21").Value = "NT"
$ws.Range("E21").Value = 6458
$ws.Range("F21").Value = "Lunglav"
$ws.Range("G21").Value = "Lobaria pulmonaria"
$ws.Range("H21").Value = "(L.) Hoffm."
$ws.Range("Q21").Value = 756485
$ws.Range("R21").Value = 7212023
$ws.Range("AJ21").Value = "sälg"
$ws.Range("AK21").Value = "Salix caprea"
$ws.Range("AO21").Value = "Salix caprea"

# Row 22
$ws.Range("A22").Value = 112171795
$ws.Range("B22").Value = 77650
$ws.Range("D22").Value = "NT"
$ws.Range("E22").Value = 6425
$ws.Range("F22").Value = "Garnlav"
$ws.Range("G22").Value = "Alectoria sarmentosa"
$ws.Range("H22").Value = "(Ach.) Ach."
$ws.Range("Q22").Value = 756378
$ws.Range("R22").Value = 7212050
$ws.Range("AJ22").ClearContents()
$ws.Range("AK22").ClearContents()
$ws.Range("AO22").ClearContents()
